$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

$ws.Cells.Item(3, 3).Value = 11.3414
$ws.Cells.Item(3, 4).Value = 20.1608
$ws.Cells.Item(3, 5).Value = 27.3292
$ws.Cells.Item(4, 2).Value = "ASALCBR"
$ws.Cells.Item(4, 3).Value = 9.1144
$ws.Cells.Item(4, 4).Value = 9.318300000000001
$ws.Cells.Item(4, 5).Value = 23.2073
$ws.Cells.Item(5, 2).Value = "INDIACEM"
$ws.Cells.Item(5, 3).Value = 8.253500000000001
$ws.Cells.Item(5, 4).Value = 8.544700000000001
$ws.Cells.Item(5, 5).Value = 10.3255
$ws.Cells.Item(6, 3).Value = 8.0824
$ws.Cells.Item(6, 4).Value = 7.0465
$ws.Cells.Item(6, 5).Value = -21.6855
$ws.Cells.Item(7, 3).Value = 7.9445
$ws.Cells.Item(7, 4).Value = 10.2264
$ws.Cells.Item(7, 5).Value = 26.4557
$ws.Cells.Item(8, 3).Value = 7.8921
$ws.Cells.Item(8, 4).Value = 5.0073
$ws.Cells.Item(8, 5).Value = -9.5098
$ws.Cells.Item(9, 3).Value = 7.3114
$ws.Cells.Item(9, 4).Value = 12.5008
$ws.Cells.Item(9, 5).Value = 14.6609
$ws.Cells.Item(10, 3).Value = 7.0443
$ws.Cells.Item(10, 4).Value = 11.3562
$ws.Cells.Item(10, 5).Value = 12.4606
$ws.Cells.Item(11, 2).Value = "PDSL"
$ws.Cells.Item(11, 3).Value = 6.4014
$ws.Cells.Item(11, 4).Value = 9.772399999999999
$ws.Cells.Item(11, 5).Value = 15.983
$ws.Cells.Item(12, 3).Value = 6.3147
$ws.Cells.Item(12, 4).Value = 12.238
$ws.Cells.Item(12, 5).Value = 14.3484
$ws.Cells.Item(13, 3).Value = 6.297
$ws.Cells.Item(13, 4).Value = 19.7672
$ws.Cells.Item(13, 5).Value = 19.8565
$ws.Cells.Item(14, 2).Value = "VENKEYS"
$ws.Cells.Item(14, 3).Value = 6.0782
$ws.Cells.Item(14, 4).Value = 6.7115
$ws.Cells.Item(14, 5).Value = 4.3904
$ws.Cells.Item(15, 3).Value = 5.3807
$ws.Cells.Item(15, 4).Value = 5.9345
$ws.Cells.Item(15, 5).Value = -0.422
$ws.Cells.Item(16, 2).Value = "MRPL"
$ws.Cells.Item(16, 3).Value = 5.3318
$ws.Cells.Item(16, 4).Value = 15.5598
$ws.Cells.Item(16, 5).Value = 26.4553
$ws.Cells.Item(17, 2).Value = "BLSE"
$ws.Cells.Item(17, 3).Value = 5.2508
$ws.Cells.Item(17, 4).Value = 4.2187
$ws.Cells.Item(17, 5).Value = -1.9651
$ws.Cells.Item(18, 3).Value = 5.0247
$ws.Cells.Item(18, 4).Value = 11.5961
$ws.Cells.Item(18, 5).Value = 8.0814
$ws.Cells.Item(21, 3).Value = 4.9235
$ws.Cells.Item(21, 4).Value = 12.3084
$ws.Cells.Item(21, 5).Value = 4.8156
$ws.Cells.Item(22, 2).Value = "VIMTALABS"
$ws.Cells.Item(22, 3).Value = 4.7909
$ws.Cells.Item(22, 4).Value = 4.9702
$ws.Cells.Item(22, 5).Value = -0.1557
$ws.Cells.Item(23, 3).Value = 4.5909
$ws.Cells.Item(23, 4).Value = 10.479
$ws.Cells.Item(23, 5).Value = 27.4792
$ws.Cells.Item(24, 2).Value = "IIFL"
$ws.Cells.Item(24, 3).Value = 4.5507
$ws.Cells.Item(24, 4).Value = 11.6374
$ws.Cells.Item(24, 5).Value = 21.0002
$ws.Cells.Item(25, 2).Value = "MARINE"
$ws.Cells.Item(25, 3).Value = 4.2865
$ws.Cells.Item(25, 4).Value = 0.8826000000000001
$ws.Cells.Item(25, 5).Value = 13.2019
$ws.Cells.Item(26, 2).Value = "DEEDEV"
$ws.Cells.Item(26, 3).Value = 4.2721
$ws.Cells.Item(26, 4).Value = -2.6649
$ws.Cells.Item(26, 5).Value = -3.4677
$ws.Cells.Item(27, 3).Value = 3.9995
$ws.Cells.Item(27, 4).Value = 11.441
$ws.Cells.Item(27, 5).Value = 9.729699999999999
$ws.Cells.Item(28, 3).Value = 3.8291
$ws.Cells.Item(28, 4).Value = 4.3335
$ws.Cells.Item(28, 5).Value = 3.8692
$ws.Cells.Item(30, 2).Value = "BLS"
$ws.Cells.Item(30, 3).Value = 3.7073
$ws.Cells.Item(30, 4).Value = 0.6441
$ws.Cells.Item(30, 5).Value = -0.6205000000000001
$ws.Cells.Item(32, 2).Value = "SKYGOLD"
$ws.Cells.Item(32, 3).Value = 3.6189
$ws.Cells.Item(32, 4).Value = -0.93
$ws.Cells.Item(32, 5).Value = 37.5878
$ws.Cells.Item(33, 2).Value = "SALASAR"
$ws.Cells.Item(33, 3).Value = 3.5533
$ws.Cells.Item(33, 4).Value = 8.5106
$ws.Cells.Item(33, 5).Value = 14.9944
$ws.Cells.Item(34, 2).Value = "SHRINGARMS"
$ws.Cells.Item(34, 3).Value = 3.5278
$ws.Cells.Item(34, 4).Value = 4.7385
$ws.Cells.Item(34, 5).Value = 24.9707
$ws.Cells.Item(35, 2).Value = "RAMASTEEL"
$ws.Cells.Item(35, 3).Value = 3.4137
$ws.Cells.Item(35, 4).Value = 3.3099
$ws.Cells.Item(35, 5).Value = 4.888
$ws.Cells.Item(36, 2).Value = "ABREL"
$ws.Cells.Item(36, 3).Value = 3.3997
$ws.Cells.Item(36, 4).Value = 11.7213
$ws.Cells.Item(36, 5).Value = 11.2584
$ws.Cells.Item(37, 2).Value = "OIL"
$ws.Cells.Item(37, 3).Value = 3.365
$ws.Cells.Item(37, 4).Value = 3.6114
$ws.Cells.Item(37, 5).Value = 5.0387
$ws.Cells.Item(38, 2).Value = "RSYSTEMS"
$ws.Cells.Item(38, 3).Value = 3.3301
$ws.Cells.Item(38, 4).Value = 4.49
$ws.Cells.Item(38, 5).Value = 6.8363
$ws.Cells.Item(39, 2).Value = "REFEX"
$ws.Cells.Item(39, 3).Value = 3.2158
$ws.Cells.Item(39, 4).Value = 0.6326000000000001
$ws.Cells.Item(39, 5).Value = 2.6224
$ws.Cells.Item(40, 2).Value = "RELTD"
$ws.Cells.Item(40, 3).Value = 3.1909
$ws.Cells.Item(40, 4).Value = 9.976800000000001
$ws.Cells.Item(40, 5).Value = -1.5443
$ws.Cells.Item(41, 3).Value = 3.16
$ws.Cells.Item(41, 4).Value = 8.67
$ws.Cells.Item(41, 5).Value = 5.7265
$ws.Cells.Item(42, 2).Value = "SAPPHIRE"
$ws.Cells.Item(42, 3).Value = 3.0365
$ws.Cells.Item(42, 4).Value = 4.8533
$ws.Cells.Item(42, 5).Value = 2.2124
$ws.Cells.Item(43, 2).Value = "INDOSTAR"
$ws.Cells.Item(43, 3).Value = 3.0133
$ws.Cells.Item(43, 4).Value = 7.8992
$ws.Cells.Item(43, 5).Value = 2.4944
$ws.Cells.Item(44, 2).Value = "MTARTECH"
$ws.Cells.Item(44, 3).Value = 3.0132
$ws.Cells.Item(44, 4).Value = 7.1148
$ws.Cells.Item(44, 5).Value = 30.9225
$ws.Cells.Item(45, 2).Value = "GPPL"
$ws.Cells.Item(45, 3).Value = 3.0058
$ws.Cells.Item(45, 4).Value = 6.5155
$ws.Cells.Item(45, 5).Value = 8.2073
$ws.Cells.Item(46, 2).Value = "MIDWESTLTD"
$ws.Cells.Item(46, 3).Value = 2.9544
$ws.Cells.Item(46, 4).Value = -1.6133
$ws.Cells.Item(46, 5).Value = "N/A"
$ws.Cells.Item(47, 2).Value = "MFSL"
$ws.Cells.Item(47, 3).Value = 2.9299
$ws.Cells.Item(47, 4).Value = 2.9842
$ws.Cells.Item(47, 5).Value = -0.7995
$ws.Cells.Item(48, 2).Value = "POWERINDIA"
$ws.Cells.Item(48, 3).Value = 2.9058
$ws.Cells.Item(48, 4).Value = 7.2941
$ws.Cells.Item(48, 5).Value = -0.0611
$ws.Cells.Item(49, 3).Value = 2.7977
$ws.Cells.Item(49, 4).Value = 3.6559
$ws.Cells.Item(49, 5).Value = -1.8397
$ws.Cells.Item(50, 2).Value = "MAMATA"
$ws.Cells.Item(50, 3).Value = 2.786
$ws.Cells.Item(50, 4).Value = 2.1957
$ws.Cells.Item(50, 5).Value = 1.3965
$ws.Cells.Item(51, 2).Value = "PSPPROJECT"
$ws.Cells.Item(51, 3).Value = 2.7716
$ws.Cells.Item(51, 4).Value = 16.8316
$ws.Cells.Item(51, 5).Value = 23.2549
$ws.Cells.Item(52, 2).Value = "DBL"
$ws.Cells.Item(52, 3).Value = 2.7038
$ws.Cells.Item(52, 4).Value = 3.7898
$ws.Cells.Item(52, 5).Value = 4.8657
$ws.Cells.Item(53, 2).Value = "CARYSIL"
$ws.Cells.Item(53, 3).Value = 2.6889
$ws.Cells.Item(53, 4).Value = 2.1731
$ws.Cells.Item(53, 5).Value = 11.063
$ws.Cells.Item(54, 2).Value = "OBEROIRLTY"
$ws.Cells.Item(54, 3).Value = 2.6847
$ws.Cells.Item(54, 4).Value = 3.5002
$ws.Cells.Item(54, 5).Value = 11.2066
$ws.Cells.Item(55, 2).Value = "CREDITACC"
$ws.Cells.Item(55, 3).Value = 2.6715
$ws.Cells.Item(55, 4).Value = 1.3146
$ws.Cells.Item(55, 5).Value = 6.5031
$ws.Cells.Item(56, 2).Value = "ACUTAAS"
$ws.Cells.Item(56, 3).Value = 2.6461
$ws.Cells.Item(56, 4).Value = 3.0119
$ws.Cells.Item(56, 5).Value = 31.5829
$ws.Cells.Item(57, 2).Value = "ALICON"
$ws.Cells.Item(57, 3).Value = 2.6145
$ws.Cells.Item(57, 4).Value = 8.6976
$ws.Cells.Item(57, 5).Value = 14.0969
$ws.Cells.Item(58, 3).Value = 2.5992
$ws.Cells.Item(58, 4).Value = 1.1439
$ws.Cells.Item(58, 5).Value = 6.7359
$ws.Cells.Item(59, 2).Value = "DBCORP"
$ws.Cells.Item(59, 3).Value = 2.5984
$ws.Cells.Item(59, 4).Value = 5.2337
$ws.Cells.Item(59, 5).Value = 1.1885
$ws.Cells.Item(60, 2).Value = "NBCC"
$ws.Cells.Item(60, 3).Value = 2.569
$ws.Cells.Item(60, 4).Value = 5.8107
$ws.Cells.Item(60, 5).Value = 10.3661
$ws.Cells.Item(61, 2).Value = "GENUSPOWER"
$ws.Cells.Item(61, 3).Value = 2.5601
$ws.Cells.Item(61, 4).Value = 10.5202
$ws.Cells.Item(61, 5).Value = 7.2408
$ws.Cells.Item(62, 2).Value = "CEATLTD"
$ws.Cells.Item(62, 3).Value = 2.5456
$ws.Cells.Item(62, 4).Value = -0.5772
$ws.Cells.Item(62, 5).Value = 21.6823
$ws.Cells.Item(63, 3).Value = 2.5442
$ws.Cells.Item(63, 4).Value = 10.5391
$ws.Cells.Item(63, 5).Value = 33.877
$ws.Cells.Item(64, 2).Value = "KERNEX"
$ws.Cells.Item(64, 3).Value = 2.5193
$ws.Cells.Item(64, 4).Value = 10.1664
$ws.Cells.Item(64, 5).Value = 30.3075
$ws.Cells.Item(65, 2).Value = "SOLEX"
$ws.Cells.Item(65, 3).Value = 2.4418
$ws.Cells.Item(65, 4).Value = 3.9749
$ws.Cells.Item(65, 5).Value = "N/A"
$ws.Cells.Item(66, 2).Value = "DIVISLAB"
$ws.Cells.Item(66, 3).Value = 2.4416
$ws.Cells.Item(66, 4).Value = 1.1831
$ws.Cells.Item(66, 5).Value = 17.2511
$ws.Cells.Item(67, 3).Value = 2.4251
$ws.Cells.Item(67, 4).Value = 1.7731
$ws.Cells.Item(67, 5).Value = 5.0657
$ws.Cells.Item(68, 2).Value = "JSFB"
$ws.Cells.Item(68, 3).Value = 2.4144
$ws.Cells.Item(68, 4).Value = 2.2765
$ws.Cells.Item(68, 5).Value = -4.6623
$ws.Cells.Item(69, 3).Value = 2.3986
$ws.Cells.Item(69, 4).Value = 2.1459
$ws.Cells.Item(69, 5).Value = 9.174300000000001
$ws.Cells.Item(70, 2).Value = "AARTIDRUGS"
$ws.Cells.Item(70, 3).Value = 2.3788
$ws.Cells.Item(70, 4).Value = 2.3788
$ws.Cells.Item(70, 5).Value = 3.9466
$ws.Cells.Item(71, 2).Value = "JAINREC"
$ws.Cells.Item(71, 3).Value = 2.3184
$ws.Cells.Item(71, 4).Value = -0.8645
$ws.Cells.Item(71, 5).Value = "N/A"
$ws.Cells.Item(72, 2).Value = "IFCI"
$ws.Cells.Item(72, 3).Value = 2.3129
$ws.Cells.Item(72, 4).Value = 4.8859
$ws.Cells.Item(72, 5).Value = 8.0296
$ws.Cells.Item(73, 2).Value = "VERANDA"
$ws.Cells.Item(73, 3).Value = 2.2876
$ws.Cells.Item(73, 4).Value = -5.2914
$ws.Cells.Item(73, 5).Value = 10.5388
$ws.Cells.Item(74, 2).Value = "ANANDRATHI"
$ws.Cells.Item(74, 3).Value = 2.2713
$ws.Cells.Item(74, 4).Value = 1.3843
$ws.Cells.Item(74, 5).Value = 11.6634
$ws.Cells.Item(76, 2).Value = "SJS"
$ws.Cells.Item(76, 3).Value = 2.2692
$ws.Cells.Item(76, 4).Value = 6.3855
$ws.Cells.Item(76, 5).Value = 10.4654
